$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expansion List")

# Existing data (row 13): 10901-7 | Display | SNOMEDCT | 2020-09 | 2.16.840.1.113883.6.96 | FN
# New rows being appended to the "Expanded Code List" table (A12:F12 header).
# Cells are populated in column order A, C, D, F, B, E to match the original
# authoring order (and resulting shared-string table order).

# Row 14: same code (10901-7) but a newer code system version (2021-09)
$ws.Range("A14").Value = "10901-7"
$ws.Range("C14").Value = "SNOMEDCT"
$ws.Range("D14").Value = "2021-09"
$ws.Range("F14").Value = "FN"
$ws.Range("B14").Value = "Display for 2021-09"
$ws.Range("E14").Value = "2.16.840.1.113883.6.96"

# Row 15: a new SNOMEDCT code (10901-8) for the 2021-09 version
$ws.Range("A15").Value = "10901-8"
$ws.Range("C15").Value = "SNOMEDCT"
$ws.Range("D15").Value = "2021-09"
$ws.Range("F15").Value = "FN"
$ws.Range("B15").Value = "Display for 10901-8"
$ws.Range("E15").Value = "2.16.840.1.113883.6.96"

# Row 16: the same code (10901-8) but under the LOINC code system
$ws.Range("A16").Value = "10901-8"
$ws.Range("C16").Value = "LOINC"
$ws.Range("D16").Value = "2021-09"
$ws.Range("F16").Value = "FN"
$ws.Range("B16").Value = "Display for 10901-8 LOINC"
$ws.Range("E16").Value = "2.16.840.1.113883.6.1"

# Match style of the existing data rows (A13:F13) for the new rows
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F16").PasteSpecial(-4122)
$ws.Range("A1").Select()
$excel.CutCopyMode = 0

# Give row 16 a slightly smaller explicit row height, matching the target file
$ws.Rows.Item(16).RowHeight = 14.25

# Update the selection/active cell to land on the newly added last row
$ws.Range("A16:XFD16").Select()
